$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used below to copy a cell's formatting onto a new cell
$xlPasteFormats = -4122

# --- Row 3: two empty, bordered cells (same style as J3/K3) ---
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial($xlPasteFormats)
$ws.Range("M3").PasteSpecial($xlPasteFormats)

# --- Row 4: year headers 2022 / 2023 (same style as K4) ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial($xlPasteFormats)
$ws.Range("M4").PasteSpecial($xlPasteFormats)
$ws.Range("L4").Value = 2022
$ws.Range("M4").Value = 2023

# --- Row 5 data (same style as K5) ---
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial($xlPasteFormats)
$ws.Range("M5").PasteSpecial($xlPasteFormats)
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 700

# --- Row 6 data (same style as K6) ---
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial($xlPasteFormats)
$ws.Range("M6").PasteSpecial($xlPasteFormats)
$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = 6

# --- Row 7 data (L7 numeric like K7, M7 dash like K6/K8) ---
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial($xlPasteFormats)
$ws.Range("K6").Copy()
$ws.Range("M7").PasteSpecial($xlPasteFormats)
$ws.Range("L7").Value = 23
$ws.Range("M7").Value = "-"

# --- Row 8 data (same style as K8) ---
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial($xlPasteFormats)
$ws.Range("M8").PasteSpecial($xlPasteFormats)
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 5

# --- Row 9 data (same style as K9) ---
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial($xlPasteFormats)
$ws.Range("M9").PasteSpecial($xlPasteFormats)
$ws.Range("L9").Value = 23
$ws.Range("M9").Value = 21

# --- Row 10 totals (same style as K10) ---
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial($xlPasteFormats)
$ws.Range("M10").PasteSpecial($xlPasteFormats)
$ws.Range("L10").Value = 172
$ws.Range("M10").Value = 143

# --- Row 11 footnote: shrink the font on the three label cells (new font/style) ---
$ws.Range("A11:C11").Font.Size = 8

# --- Page setup: A4-ish paper (paperSize 256), portrait orientation ---
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1
